$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 246, shifting existing rows 246:343 down to 247:344
$ws.Rows(246).Insert()

# Fill in the values for the newly inserted row 246
$ws.Range("A246").Value = 10
$ws.Range("B246").Value = "Vega Modelo de Temuco"
$ws.Range("C246").Value = "La Araucanía"
$ws.Range("D246").Value = 44755
$ws.Range("E246").Value = 9
$ws.Range("F246").Value = 100112009
$ws.Range("G246").Value = "Acelga"
$ws.Range("H246").Value = "Sin especificar"
$ws.Range("I246").Value = "Primera"
$ws.Range("J246").Value = 50
$ws.Range("K246").Value = 8000
$ws.Range("L246").Value = 8000
$ws.Range("M246").Value = 8000
$ws.Range("N246").Value = "$/docena de atados (12 kilos)"
$ws.Range("O246").Value = "Región Metropolitana"
$ws.Range("P246").Value = 667
$ws.Range("Q246").Value = 12
$ws.Range("R246").Value = "Hortaliza"
